# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 2
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 1
    23 = 0
    24 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
